$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.InsertAfter("`r`rSimplify slides, less graphs, explain algorithms with diagrams not words")
